# The dataset gained one new weekly record. A new row is inserted at row 25
# (pushing the existing rows 25-68 down to 26-69, which also grows the used
# range from A1:R68 to A1:R69), and the newly inserted row 25 is populated
# with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = 9
$ws.Range("B25").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44477
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 100112005
$ws.Range("G25").Value = "Puerro"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 7500
$ws.Range("N25").Value = "$/paquete 20 unidades"
$ws.Range("O25").Value = "Provincia de Chacabuco"
$ws.Range("P25").Value = 375
$ws.Range("Q25").Value = 20
$ws.Range("R25").Value = "Hortaliza"
